# Weekly update: insert the latest week's price record at the top of the
# date-ordered data block (new row 11), pushing all existing records
# (old rows 11-28) down by one row (new rows 12-29).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before the current row 11; this shifts rows 11-28
# down to 12-29 and extends the used range / dimension accordingly.
$ws.Rows.Item(11).Insert()

# Populate the new row 11 with this week's data.
$ws.Range("A11").Value = 2
$ws.Range("B11").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C11").Value = "Coquimbo"
$ws.Range("D11").Value = 44482
$ws.Range("E11").Value = 4
$ws.Range("F11").Value = 100112032
$ws.Range("G11").Value = "Zapallo italiano"
$ws.Range("H11").Value = "Sin especificar"
$ws.Range("I11").Value = "Primera"
$ws.Range("J11").Value = 400
$ws.Range("K11").Value = 11000
$ws.Range("L11").Value = 12000
$ws.Range("M11").Value = 11500
$ws.Range("N11").Value = "$/caja 60 unidades"
$ws.Range("O11").Value = "Provincia de Limarí"
$ws.Range("P11").Value = 192
$ws.Range("Q11").Value = 60
$ws.Range("R11").Value = "Hortaliza"
